$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "identificado con la clave IYE",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "identificado con la clave ",
    2
)
